# Insert a new "VENCIMENTO" column between the existing TOTAL (B) and
# TIPO (C) columns, shifting TIPO to column D, and populate it with the
# due-day values. Also apply currency formatting to the TOTAL column and
# a centered integer format to the new VENCIMENTO column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift column C (TIPO) to the right to make room for the new column.
$ws.Columns.Item(3).Insert()

# Header
$ws.Range("C1").Value = "VENCIMENTO"

# Due-day values for rows 2-8
$dias = @(26, 15, 15, 10, 16, 10, 25)
for ($i = 0; $i -lt $dias.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $dias[$i]
}

# Number formatting: apply the named "Currency" (Moeda) cell style, then
# refine the format mask to the Brazilian Real pattern used by the workbook.
$ws.Range("B1:B8").Style = "Currency"
$ws.Range("B1:B8").NumberFormat = '_-"R$"\ * #,##0.00_-;\-"R$"\ * #,##0.00_-;_-"R$"\ * "-"??_-;_-@_-'
$ws.Range("C1:C8").NumberFormat = "0"
$ws.Range("C1:C8").HorizontalAlignment = -4108

# Column widths (best-fit to content, mirroring Excel's AutoFit result)
$ws.Columns.Item(2).ColumnWidth = 11.3
$ws.Columns.Item(3).ColumnWidth = 12.3

# Selection as in the target file
$ws.Range("F5").Select()
